$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "Linkdin" header (G1) to the corrected spelling "linkedin"
$ws.Range("G1").Value = "linkedin"

# Update the remembered selection to match the authored state (H10)
$ws.Range("H10").Select()
